$d = $word.ActiveDocument

# 1. "Niveau: Brugermål" -> "Niveau: Underfunktion"
$d.Content.Find.Execute(": Brugermål", $false, $false, $false, $false, $false, $true, 1, $false, ": Underfunktion", 2)

# 2. Merge split runs "2a. Der kommer ikke noget svar / en fejl fra " + "banken"
$d.Content.Find.Execute("2a. Der kommer ikke noget svar / en fejl fra banken", $false, $false, $false, $false, $false, $true, 1, $false, "2a. Der kommer ikke noget svar / en fejl fra banken", 2)

# 3. Merge split runs "Systemet samarbejder med " + "banken" + "s "
$d.Content.Find.Execute("Systemet samarbejder med bankens ", $false, $false, $false, $false, $false, $true, 1, $false, "Systemet samarbejder med bankens ", 2)

# 4. Merge split runs "Brugerinterfacet må ikke blive påvirket af at der afventes svar fra " + "banken"
$d.Content.Find.Execute("Brugerinterfacet må ikke blive påvirket af at der afventes svar fra banken", $false, $false, $false, $false, $false, $true, 1, $false, "Brugerinterfacet må ikke blive påvirket af at der afventes svar fra banken", 2)
